$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Cells whose new value is a plain number-looking string (e.g. "310.73").
# Excel normally auto-converts such literals typed into a Range.Value to a
# numeric cell; the source data must stay TEXT (it is a scraped price column
# that elsewhere contains things like "26.873.65"). Force text by switching the
# cell to the Text number format before the assignment, then restore the
# default "Normal" style so no stray formatting is left behind.
$numericLooking = @("D5", "D7", "D8", "D9", "D10", "D11", "D14", "D15", "D16", "D21", "D22", "D24", "D25", "D27", "D28", "D29", "D30", "D31", "D32", "D33", "D34", "D35", "D36", "D37", "D38", "D39", "D40", "D41", "D43", "D44", "D46", "D47", "D48", "D49", "D51")
foreach ($ref in $numericLooking) {
    $ws.Range($ref).NumberFormat = "@"
}

# --- Apply the updated values cell by cell (prices, 1h-volume deltas, and
# the Hedera / MXToken row swap with their refreshed figures).
$ws.Range("D2").Value = '26.873.71'
$ws.Range("E2").Value = '  -1.77%  '
$ws.Range("D3").Value = '1.825.13'
$ws.Range("E3").Value = '  -1.73%  '
$ws.Range("E4").Value = '  +0.58%  '
$ws.Range("D5").Value = '310.73'
$ws.Range("E5").Value = '  -1.14%  '
$ws.Range("D7").Value = '0.4580'
$ws.Range("E7").Value = '  -0.76%  '
$ws.Range("D8").Value = '0.3672'
$ws.Range("E8").Value = '  -1.16%  '
$ws.Range("D9").Value = '0.07150'
$ws.Range("E9").Value = '  -2.34%  '
$ws.Range("D10").Value = '0.8712'
$ws.Range("E10").Value = '  -1.20%  '
$ws.Range("D11").Value = '0.07784'
$ws.Range("E11").Value = '  -0.30%  '
$ws.Range("E12").Value = '  -2.03%  '
$ws.Range("D13").Value = '1.819.61'
$ws.Range("E13").Value = '  -2.99%  '
$ws.Range("D14").Value = '5.315'
$ws.Range("E14").Value = '  -1.26%  '
$ws.Range("D15").Value = '6.373'
$ws.Range("E15").Value = '  -2.67%  '
$ws.Range("D16").Value = '86.76'
$ws.Range("E16").Value = '  -5.57%  '
$ws.Range("E18").Value = '  -4.17%  '
$ws.Range("E19").Value = '  +0.49%  '
$ws.Range("D20").Value = '26.896.23'
$ws.Range("E20").Value = '  -1.73%  '
$ws.Range("D21").Value = '14.43'
$ws.Range("E21").Value = '  -2.49%  '
$ws.Range("D22").Value = '4.987'
$ws.Range("E22").Value = '  -2.79%  '
$ws.Range("E23").Value = '  -0.63%  '
$ws.Range("D24").Value = '2.006'
$ws.Range("E24").Value = '  +3.77%  '
$ws.Range("D25").Value = '150.98'
$ws.Range("E25").Value = '  -0.86%  '
$ws.Range("E26").Value = '  -1.03%  '
$ws.Range("D27").Value = '1.956'
$ws.Range("E27").Value = '  -5.65%  '
$ws.Range("D28").Value = '113.44'
$ws.Range("E28").Value = '  -2.22%  '
$ws.Range("D29").Value = '4.913'
$ws.Range("E29").Value = '  -3.81%  '
$ws.Range("D30").Value = '0.08799'
$ws.Range("E30").Value = '  -0.69%  '
$ws.Range("D31").Value = '3.012'
$ws.Range("E31").Value = '  -0.97%  '
$ws.Range("D32").Value = '0.7441'
$ws.Range("E32").Value = '  -3.68%  '
$ws.Range("D33").Value = '4.469'
$ws.Range("E33").Value = '  -0.60%  '
$ws.Range("D34").Value = '1.127'
$ws.Range("E34").Value = '  -4.14%  '
$ws.Range("D35").Value = '2.506'
$ws.Range("E35").Value = '  -5.47%  '
$ws.Range("D36").Value = '1.086'
$ws.Range("E36").Value = '  +0.82%  '
$ws.Range("D37").Value = '0.01933'
$ws.Range("E37").Value = '  -1.28%  '
$ws.Range("B38").Value = 'Hedera'
$ws.Range("C38").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D38").Value = '0.05108'
$ws.Range("E38").Value = '  -2.34%  '
$ws.Range("B39").Value = 'MXToken'
$ws.Range("C39").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D39").Value = '2.902'
$ws.Range("E39").Value = '  -1.69%  '
$ws.Range("D40").Value = '6.908'
$ws.Range("E40").Value = '  -1.84%  '
$ws.Range("D41").Value = '0.4953'
$ws.Range("E41").Value = '  -3.68%  '
$ws.Range("E42").Value = '  -2.65%  '
$ws.Range("D43").Value = '8.264'
$ws.Range("E43").Value = '  -1.93%  '
$ws.Range("D44").Value = '0.4665'
$ws.Range("E44").Value = '  -3.26%  '
$ws.Range("E45").Value = '  +0.58%  '
$ws.Range("D46").Value = '10.04'
$ws.Range("E46").Value = '  -2.90%  '
$ws.Range("D47").Value = '101.30'
$ws.Range("E47").Value = '  -1.82%  '
$ws.Range("D48").Value = '1.603'
$ws.Range("E48").Value = '  -2.92%  '
$ws.Range("D49").Value = '0.06082'
$ws.Range("E49").Value = '  -2.28%  '
$ws.Range("E50").Value = '  -2.44%  '
$ws.Range("D51").Value = '36.62'

# --- Drop the temporary Text format back to Normal so the saved styles match
# the original (no explicit cell style on the data rows).
foreach ($ref in $numericLooking) {
    $ws.Range($ref).Style = "Normal"
}
